$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 (ck1.txt row): value 39 -> 70
$ws.Range("C2").Value() = 70

# Remove row 3 (the ck2.txt entry) entirely; this shifts the ck3.txt row (old row 4,
# which carries its own custom cell style) up to become row 3
$ws.Rows("3").Delete() | Out-Null

# Carry the ck3.txt row's cell formatting over onto the updated C2 cell, matching
# the other "70" cell's look
$ws.Range("C3").Select() | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C2").Select() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null

# Add a new, empty, underlined cell at D5
$ws.Range("D5").Font.Underline = $true

# Move the active selection to D5, matching the saved view state
$ws.Range("D5").Select() | Out-Null
